$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newline = [char]10

$ws.Range("B20").Value = "11-50" + $newline + "(up to 100 for construction sector)"
$ws.Range("B21").Value = "51-250" + $newline + "(up to 400 for construction sector)"
$ws.Range("B22").Value = ">250" + $newline + "(> 400 for construction sector)"
